# UiComponentClassDiagram.pptx - "Update dev guide until 2.2 UI component" edit.
#
# The UI class diagram on slide 1 renames two class boxes that used to model a
# "Person" domain object so that they model a "Task" instead:
#   "PersonListPanel" -> "Task" + "ListPanel"   (typed as two runs, same as PowerPoint
#                                                 does when a user edits in place)
#   "PersonCard"       -> "Task" + "Card"
#
# We find the shapes by their current text (robust to any shape-index churn)
# rather than hard-coding a shape index, then retype the text the way a user
# would: replace the whole run with "Task", then continue typing the
# remaining suffix right after it. PowerPoint's editor keeps that as two
# separate <a:r> runs (one per "typing burst"), which is exactly what the
# target XML shows.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# NOTE: this PowerShell host only binds positional parameters, so the
# helper below is called positionally (Name -Param value style binding
# silently fails to bind).
function Set-ClassBoxText {
    param([string]$OldText, [string]$NewPrefix, [string]$NewSuffix)

    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $OldText) {
                $tr.Text = $NewPrefix
                $tr.InsertAfter($NewSuffix)
                return
            }
        }
    }

    throw ("Shape with text '" + $OldText + "' not found")
}

Set-ClassBoxText "PersonListPanel" "Task" "ListPanel"
Set-ClassBoxText "PersonCard" "Task" "Card"
